# -----------------------------------------------------------------------
# Add a "Feedback" section (heading + contact paragraph with a mailto
# hyperlink) at the very end of the document, replacing the old
# "Last updated: April 16th, 2020." line.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the paragraph that currently reads
#    "Last updated: April 16th, 2020." - this is the last paragraph
#    of the document, right after the contacts table.
# ------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.ClearFormatting()
$found = $findRng.Find.Execute("Last updated: April 16th, 2020.")
if (-not $found) {
    Write-Host "ERROR: could not locate the 'Last updated' paragraph"
}

# ------------------------------------------------------------------
# 2. Insert a new paragraph ("Feedback") immediately before it. This
#    reuses the existing trailing paragraph mark of the document
#    (Word will not allow deleting the very last paragraph mark), so
#    instead of deleting a paragraph we push a new one in front of it
#    and then turn the old paragraph into the feedback text paragraph.
# ------------------------------------------------------------------
$insertPoint = $findRng.Duplicate
$insertPoint.Collapse(1) | Out-Null   # wdCollapseStart
$insertPoint.InsertBefore("Feedback`r")

# Style the new "Feedback" paragraph as Heading 2.
$headingRng = $d.Content
$headingRng.Find.ClearFormatting()
$headingRng.Find.Execute("Feedback") | Out-Null
$headingStart = $headingRng.Start
$headingRng.Paragraphs.Item(1).Range.Style = "Heading 2"

# ------------------------------------------------------------------
# 3. Replace the text of the old "Last updated..." paragraph (which
#    is now the paragraph right after the heading) with the feedback
#    sentence, and restyle it as Block Text.
# ------------------------------------------------------------------
$replaceRng = $d.Content
$replaceRng.Find.ClearFormatting()
$replaceRng.Find.Execute( `
    "Last updated: April 16th, 2020.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "If you have any questions or comments about this guidance, such as suggestions for improvements, please contact: itpolicycontent@digital.justice.gov.uk.", `
    2) | Out-Null

$bodyRng = $d.Content
$bodyRng.Find.ClearFormatting()
$bodyRng.Find.Execute("If you have any questions") | Out-Null
$bodyRng.Paragraphs.Item(1).Range.Style = "Block Text"

# ------------------------------------------------------------------
# 4. Turn the e-mail address into a real mailto: hyperlink.
# ------------------------------------------------------------------
$linkRng = $d.Content
$linkRng.Find.ClearFormatting()
$linkRng.Find.Execute("itpolicycontent@digital.justice.gov.uk") | Out-Null
$d.Hyperlinks.Add($linkRng, "mailto:itpolicycontent@digital.justice.gov.uk") | Out-Null

# ------------------------------------------------------------------
# 5. Wrap the new "Feedback" heading + paragraph in their own
#    "ariaid-title14" bookmark, matching the pattern used for every
#    other top-level section in this document.
# ------------------------------------------------------------------
$d.Bookmarks.Add("ariaid-title14", $d.Range($headingStart, $d.Content.End)) | Out-Null

Write-Host "Feedback section added successfully."
